$d = $word.ActiveDocument

# "Pipe" summary table (Pipe | Fixture Units | Flow Rate (gpm) | Flow Rate (m3/s) | Comment).
# Row for pipe P11 (row 12, header is row 1) was left blank; fill in the
# same values used for the other small-fixture pipes (P9/P10).
$pipeTable = $d.Tables.Item(2)
$pipeTable.Cell(12, 2).Range.Text = "2"
$pipeTable.Cell(12, 3).Range.Text = "2"
$pipeTable.Cell(12, 4).Range.Text = "0.000126"

# "Pipe diameters" table (Pipe | Flow Rate (gpm) | Diameter (in) | Diameter (m) | Justification Comments).
# Same P11 row was missing its flow rate figure.
$diameterTable = $d.Tables.Item(3)
$diameterTable.Cell(12, 2).Range.Text = "2"
